# Fix Metric Data, Allow Custom Data Upload
#
# 1. Replace the short "type" category codes on Sheet1 with descriptive
#    labels.
# 2. Insert a new timeline entry (Unified Support Teams developed /
#    Central Valley Task Force) as row 75, shifting all later rows down.
# 3. Add the new "TR" source/url lookup row to Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- 1. Re-label the "type" column values on Sheet1 --------------------
$typeMap = @{
    "epi"          = "Epidemiologic Milestone"
    "declarations" = "Declarations and Announcementss"
    "npi"          = "Nonpharmaceutical Intervention"
    "vax"          = "Vaccination"
    "test"         = "Testing"
    "treat"        = "Medications"
}

$lastRow = $ws1.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws1.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($typeMap.ContainsKey($val)) {
        $cell.Value = $typeMap[$val]
    }
}

# --- 2. Insert the new row at position 75 -------------------------------
$ws1.Rows.Item(75).Insert()

$ws1.Range("A75").Value = [DateTime]"2020-07-28"
$ws1.Range("B75").Value = "Unified Support Teams develpoed"
$ws1.Range("C75").Value = "Declarations and Announcementss"
$ws1.Range("D75").Value = "TR"
$ws1.Range("E75").Value = 6
$ws1.Range("F75").Value = "Governor Gavin Newsom announced a call to action to slow the spread of COVID-19 in these hard-hit communities. On Tuesday, July 28, 2020, the California Governor’s Office of Emergency Services (Cal OES) and the California Health and Human Services Agency (CHHS) responded by establishing the Central Valley Task Force."

# --- 3. Add the new source/url lookup row on Sheet2 ---------------------
$ws2.Range("A11").Value = "TR"
$ws2.Range("B11").Value = "Timeline provided in private communicaton from Raymundo, Trudy@CDPH"

# --- Cosmetic: leave the selection near where the new data landed -------
$ws1.Activate()
$ws1.Range("G148").Select()
$ws2.Range("G13").Select()
$ws1.Activate()
